# The commit swaps the presentation's theme palette back from the custom
# "Integral" design (green/teal accents) to the stock "Office Theme" palette
# (blue/orange accents) -- i.e. ppt/theme/theme2.xml (the theme actually
# wired to the slide master via slideMaster1.xml.rels) goes from the
# Integral clrScheme to the Office clrScheme. Font scheme / format scheme
# are identical between the two themes already, so only the 12 theme colors
# need to change.

$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

function BGRInt($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the stock "Office Theme" clrScheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) -- MsoThemeColorSchemeIndex order.
$tcs.Item(1).RGB  = BGRInt 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = BGRInt 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = BGRInt 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = BGRInt 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = BGRInt 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = BGRInt 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = BGRInt 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = BGRInt 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = BGRInt 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = BGRInt 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = BGRInt 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = BGRInt 0x95 0x4F 0x72   # folHlink
